$wb = $excel.ActiveWorkbook

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 60461.39
$ws.Cells.Item(62, 9).Value = 92991.37
$ws.Cells.Item(62, 10).Value = 9342.857
$ws.Cells.Item(62, 11).Value = 92991.37
$ws.Cells.Item(62, 12).Value = 9342.857
$ws.Cells.Item(62, 13).Value = -92367.37
$ws.Cells.Item(62, 14).Value = -10590.857

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 60461.39
$ws.Cells.Item(65, 9).Value = 92991.37
$ws.Cells.Item(65, 10).Value = 9342.857
$ws.Cells.Item(65, 11).Value = 464956.85
$ws.Cells.Item(65, 12).Value = 46714.285
$ws.Cells.Item(65, 13).Value = -461836.85
$ws.Cells.Item(65, 14).Value = -52954.285

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 2825.7778
$ws.Cells.Item(125, 9).Value = 1510.6666
$ws.Cells.Item(125, 10).Value = 3483.3333
$ws.Cells.Item(125, 11).Value = 13595.9994
$ws.Cells.Item(125, 12).Value = 31349.9997
$ws.Cells.Item(125, 13).Value = -11135.9994
$ws.Cells.Item(125, 14).Value = -36269.9997

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 45462916
$ws.Cells.Item(132, 9).Value = 55564628
$ws.Cells.Item(132, 10).Value = 5215.6665
$ws.Cells.Item(132, 11).Value = 166693884
$ws.Cells.Item(132, 12).Value = 15646.9995
$ws.Cells.Item(132, 13).Value = -166691354
$ws.Cells.Item(132, 14).Value = -20706.9995

# ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(140, 8).Value = 37765.57
$ws.Cells.Item(140, 9).Value = 34400
$ws.Cells.Item(140, 10).Value = 38326.5
$ws.Cells.Item(140, 11).Value = 34400
$ws.Cells.Item(140, 12).Value = 38326.5
$ws.Cells.Item(140, 13).Value = -29220
$ws.Cells.Item(140, 14).Value = -48686.5

# ARM row 82
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(82, 8).Value = 48630
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 48630
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 12).Value = 48630
$ws.Cells.Item(82, 14).Value = -49352

# ARM row 85
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(85, 8).Value = 48630
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 48630
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 48630
$ws.Cells.Item(85, 14).Value = -51126

# ARM row 95
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(95, 8).Value = 41000
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 41000
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 12).Value = 41000
$ws.Cells.Item(95, 14).Value = -46492

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2942796.2
$ws.Cells.Item(132, 9).Value = 894.4286
$ws.Cells.Item(132, 10).Value = 9807234
$ws.Cells.Item(132, 11).Value = 2683.2858
$ws.Cells.Item(132, 12).Value = 29421702
$ws.Cells.Item(132, 13).Value = -153.2857999999997
$ws.Cells.Item(132, 14).Value = -29426762

# CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(51, 8).Value = 20000
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 20000
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 20000
$ws.Cells.Item(51, 14).Value = -21472

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 28572066
$ws.Cells.Item(58, 9).Value = 50000640
$ws.Cells.Item(58, 10).Value = 635.73334
$ws.Cells.Item(58, 11).Value = 50000640
$ws.Cells.Item(58, 12).Value = 635.73334
$ws.Cells.Item(58, 13).Value = -50000437
$ws.Cells.Item(58, 14).Value = -1041.73334

# CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(61, 8).Value = 20000
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 20000
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 20000
$ws.Cells.Item(61, 14).Value = -20696

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 10002.546
$ws.Cells.Item(94, 9).Value = 100000
$ws.Cells.Item(94, 10).Value = 1002.8
$ws.Cells.Item(94, 11).Value = 100000
$ws.Cells.Item(94, 12).Value = 1002.8
$ws.Cells.Item(94, 13).Value = -99549
$ws.Cells.Item(94, 14).Value = -1904.8

# CRP row 109
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(109, 8).Value = 23296.25
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = 23296.25
$ws.Cells.Item(109, 11).Value = 0
$ws.Cells.Item(109, 12).Value = 23296.25
$ws.Cells.Item(109, 14).Value = -25376.25

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 28572066
$ws.Cells.Item(136, 9).Value = 50000640
$ws.Cells.Item(136, 10).Value = 635.73334
$ws.Cells.Item(136, 11).Value = 150001920
$ws.Cells.Item(136, 12).Value = 1907.20002
$ws.Cells.Item(136, 13).Value = -149999370
$ws.Cells.Item(136, 14).Value = -7007.20002

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 35357150
$ws.Cells.Item(5, 9).Value = 48611530
$ws.Cells.Item(5, 10).Value = 12155.556
$ws.Cells.Item(5, 11).Value = 145834590
$ws.Cells.Item(5, 12).Value = 36466.66800000001
$ws.Cells.Item(5, 13).Value = -145834478
$ws.Cells.Item(5, 14).Value = -36690.66800000001

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 93.666664
$ws.Cells.Item(23, 9).Value = 30
$ws.Cells.Item(23, 10).Value = 106.4
$ws.Cells.Item(23, 11).Value = 90
$ws.Cells.Item(23, 12).Value = 319.2
$ws.Cells.Item(23, 13).Value = 145
$ws.Cells.Item(23, 14).Value = -789.2

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 572.3333
$ws.Cells.Item(117, 9).Value = 286.8
$ws.Cells.Item(117, 10).Value = 2000
$ws.Cells.Item(117, 11).Value = 860.4000000000001
$ws.Cells.Item(117, 12).Value = 6000
$ws.Cells.Item(117, 13).Value = 2581.6
$ws.Cells.Item(117, 14).Value = -12884

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 932.78
$ws.Cells.Item(131, 9).Value = 464.875
$ws.Cells.Item(131, 10).Value = 973.4674
$ws.Cells.Item(131, 11).Value = 1394.625
$ws.Cells.Item(131, 12).Value = 2920.4022
$ws.Cells.Item(131, 13).Value = 3645.375
$ws.Cells.Item(131, 14).Value = -13000.4022

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 35357150
$ws.Cells.Item(135, 9).Value = 48611530
$ws.Cells.Item(135, 10).Value = 12155.556
$ws.Cells.Item(135, 11).Value = 437503770
$ws.Cells.Item(135, 12).Value = 109400.004
$ws.Cells.Item(135, 13).Value = -437501235
$ws.Cells.Item(135, 14).Value = -114470.004

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10449.117
$ws.Cells.Item(70, 9).Value = 13633.272
$ws.Cells.Item(70, 10).Value = 4611.5
$ws.Cells.Item(70, 11).Value = 13633.272
$ws.Cells.Item(70, 12).Value = 4611.5
$ws.Cells.Item(70, 13).Value = -13363.272
$ws.Cells.Item(70, 14).Value = -5151.5

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 10449.117
$ws.Cells.Item(73, 9).Value = 13633.272
$ws.Cells.Item(73, 10).Value = 4611.5
$ws.Cells.Item(73, 11).Value = 13633.272
$ws.Cells.Item(73, 12).Value = 4611.5
$ws.Cells.Item(73, 13).Value = -12697.272
$ws.Cells.Item(73, 14).Value = -6483.5

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 22732788
$ws.Cells.Item(122, 9).Value = 38470124
$ws.Cells.Item(122, 10).Value = 1083.6666
$ws.Cells.Item(122, 11).Value = 115410372
$ws.Cells.Item(122, 12).Value = 3250.9998
$ws.Cells.Item(122, 13).Value = -115407922
$ws.Cells.Item(122, 14).Value = -8150.9998

# LTW row 81
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(81, 8).Value = 20140.5
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 20140.5
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 12).Value = 20140.5
$ws.Cells.Item(81, 14).Value = -22136.5

# LTW row 84
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(84, 8).Value = 20140.5
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 20140.5
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 60421.5
$ws.Cells.Item(84, 14).Value = -70405.5

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 13789.223
$ws.Cells.Item(122, 9).Value = 16586.857
$ws.Cells.Item(122, 10).Value = 3997.5
$ws.Cells.Item(122, 11).Value = 49760.571
$ws.Cells.Item(122, 12).Value = 11992.5
$ws.Cells.Item(122, 13).Value = -47310.571
$ws.Cells.Item(122, 14).Value = -16892.5

# WVR row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 30000
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 30000
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 30000
$ws.Cells.Item(45, 14).Value = -30982

# WVR row 108
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(108, 8).Value = 37295.332
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 37295.332
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).Value = 37295.332
$ws.Cells.Item(108, 14).Value = -44975.332
